$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old trailing headers (Beleska / Uporabnik / Vreme) first so their
# shared strings get garbage-collected before we introduce the new ones.
$ws.Range("F1:H1").Clear()

# Rename surviving headers / retarget the others. Order matters: the engine
# appends newly-introduced shared strings at the end of the table, in the
# order they're first written, and garbage-collects strings that drop to
# zero references. Writing D1, then A1, then E1 (then the CCNA cells)
# reproduces the exact target shared-string order:
#   Pozicija, Ident, Kolicina, St. Dokumenta, Ekstra, CCNA
$ws.Range("D1").Value = "Količina"
$ws.Range("A1").Value = "Št. Dokumenta"
$ws.Range("E1").Value = "Ekstra"
# B1 ("Pozicija") and C1 ("Ident") keep their original text untouched.

# Column A: key formula. A2 is a standalone formula; A3:A15 share one formula.
$ws.Range("A2").Formula = "=2301150000001"
$ws.Range("A3:A15").Formula = "=2301150000001"

# Column B: running position counter. B2 is a literal 1; B3 is a standalone
# "=B2+1"; B4:B15 share one "=B3+1"-style formula.
$ws.Range("B2").Value = 1
$ws.Range("B3").Formula = "=B2+1"
$ws.Range("B4:B15").Formula = "=B3+1"

# Columns C/D/E: Ident, Kolicina (quantity) and the new "CCNA" note column.
$idents = @(104118, 104134, 104207, 104212, 104305, 106810, 117410, 120006, 251400, 258310, 258501, 259908, 259909, 260001)
$qty    = @(77000000, 10000000, 14000000, 12000000, 13000000, 418000000, 76000000, 15000000, 11000000, 11000000, 11000000, 10000000, 18000000, 60000000)

for ($i = 0; $i -lt $idents.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 3).Value = $idents[$i]
    $ws.Cells.Item($r, 4).Value = $qty[$i]
    $ws.Cells.Item($r, 4).NumberFormat = "#,##0"
    $ws.Cells.Item($r, 5).Value = "CCNA"
}

# Column widths.
$ws.Columns.Item(1).ColumnWidth = 16.14
$ws.Columns.Item(4).ColumnWidth = 11.43

# Page setup (portrait, A4 paper size code 9).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Match the saved selection/active cell from the edited workbook.
$ws.Range("I11").Select()
